$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate B2 with the new entry "Apple Limitada" instead of the previous
# duplicated "Apple Inc" value, so all entries from the source data are
# correctly reflected.
$ws.Range("B2").Value = "Apple Limitada"

# Update the active selection to reflect where the user left off editing.
$ws.Range("B2").Select()
